$wb = $excel.ActiveWorkbook

$wsBuses = $wb.Worksheets.Item("buses")
$wsLines = $wb.Worksheets.Item("lines")

# ----- buses sheet value updates -----
$wsBuses.Range("B3").Value = -0.1
$wsBuses.Range("B4").Value = -0.1

# ----- lines sheet value updates -----
$wsLines.Range("C2").Value = 0.01
$wsLines.Range("D2").Value = 0.01
$wsLines.Range("C3").Value = 0.01
$wsLines.Range("D3").Value = 0.01
$wsLines.Range("D4").Value = 0.01
$wsLines.Range("D5").Value = 0.01

# D6:D12 were formulas referencing the previous row (shared formula group
# anchored at D6). Replace them with plain literal values so the column
# becomes constant, which breaks them out of the D:H shared formula group.
$wsLines.Range("D6").Value = 0.01
$wsLines.Range("D7").Value = 0.01
$wsLines.Range("D8").Value = 0.01
$wsLines.Range("D9").Value = 0.01
$wsLines.Range("D10").Value = 0.01
$wsLines.Range("D11").Value = 0.01
$wsLines.Range("D12").Value = 0.01

# Re-enter the E:H formulas so Excel re-derives the shared formula group,
# now anchored at E6 (since D6 no longer participates).
$wsLines.Range("E6:H12").Formula = "=E5"

# D12 previously carried a distinct cell style (applyFont); clear it back
# to the default "Normal" style now that it's a plain value like the rest
# of the column.
$wsLines.Range("D12").Style = "Normal"

# ----- restore the selections recorded in the saved file -----
$wsBuses.Range("B3").Select()
$wsLines.Range("D5").Select()
